$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '54.258.99'
$ws.Range("E2").Value = '  +0.28%  '
$ws.Range("D3").Value = '2.264.48'
$ws.Range("E3").Value = '  -1.04%  '
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '495.82'
$ws.Range("E5").Value = '  -0.04%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '128.81'
$ws.Range("E6").Value = '  +0.13%  '
$ws.Range("E7").Value = '  -0.16%  '
$ws.Range("E8").Value = '  -0.95%  '
$ws.Range("E9").Value = '  +0.41%  '
$ws.Range("E10").Value = '  +0.88%  '
$ws.Range("E11").Value = '  +2.78%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.82'
$ws.Range("E12").Value = '  +3.78%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.94'
$ws.Range("E13").Value = '  +5.16%  '
$ws.Range("D14").Value = '2.664.58'
$ws.Range("E14").Value = '  -1.02%  '
$ws.Range("D15").Value = '54.228.31'
$ws.Range("E15").Value = '  +0.07%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000130'
$ws.Range("E16").Value = '  +0.12%  '
$ws.Range("D17").Value = '2.269.15'
$ws.Range("E17").Value = '  -0.71%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '10.23'
$ws.Range("E18").Value = '  +1.76%  '
$ws.Range("E19").Value = '  +0.29%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '303.00'
$ws.Range("E20").Value = '  +0.62%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.33'
$ws.Range("E21").Value = '  -1.99%  '
$ws.Range("E22").Value = '  -0.14%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '60.71'
$ws.Range("E23").Value = '  -3.18%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.998'
$ws.Range("E24").Value = '  -1.40%  '
$ws.Range("E25").Value = '  +0.29%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.32'
$ws.Range("E26").Value = '  +3.52%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '170.76'
$ws.Range("E27").Value = '  +0.99%  '
$ws.Range("E28").Value = '  -0.18%  '
$ws.Range("E29").Value = '  +1.39%  '
$ws.Range("D30").Value = '0.0₃0690'
$ws.Range("E30").Value = '  +0.05%  '
$ws.Range("E31").Value = '  +0.76%  '
$ws.Range("E32").Value = '  -0.02%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '17.78'
$ws.Range("E33").Value = '  +0.50%  '
$ws.Range("E34").Value = '  +0.09%  '
$ws.Range("E35").Value = '  +3.64%  '
$ws.Range("E36").Value = '  +0.10%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.71'
$ws.Range("E38").Value = '  -0.12%  '
$ws.Range("E39").Value = '  -0.53%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.37'
$ws.Range("E40").Value = '  +0.12%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.82'
$ws.Range("E41").Value = '  +0.67%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '124.72'
$ws.Range("E42").Value = '  -2.01%  '
$ws.Range("E43").Value = '  +1.79%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0894'
$ws.Range("E44").Value = '  +0.74%  '
$ws.Range("E45").Value = '  -0.24%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '241.55'
$ws.Range("E46").Value = '  +1.18%  '
$ws.Range("E48").Value = '  +0.82%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '10.81'
$ws.Range("E49").Value = '  +0.80%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '16.12'
$ws.Range("E50").Value = '  -1.32%  '
$ws.Range("E51").Value = '  -0.41%  '
